# Weekly update: insert two new price-report rows for "Perejil" (Vega Modelo
# de Temuco) and shift the existing data down accordingly.
#
# Row 148 gets a brand-new record (date serial 44567 = 2022-01-06), pushing
# every existing row from 148 down by one.
# A second brand-new record (date serial 44568 = 2022-01-07) is inserted
# right before the row that used to be the old row 238 (now sitting at row
# 239 after the first shift), pushing it and everything after it down by one
# more - matching the final dimension A1:R248.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToAdd = @(
    @{Row=148; Fecha=44567; Volumen=30; PrecioMin=5000; PrecioMax=5000; PrecioProm=5000; Unidad="`$/docena de atados (3 kilos)"; Origen="Provincia de Cautín"; PrecioKg=1667},
    @{Row=239; Fecha=44568; Volumen=30; PrecioMin=5000; PrecioMax=5000; PrecioProm=5000; Unidad="`$/docena de atados (3 kilos)"; Origen="Provincia de Cautín"; PrecioKg=1667}
)

foreach ($item in $rowsToAdd) {
    $r = $item.Row
    $ws.Rows.Item($r).Insert()

    $ws.Range("A$r").Value = 10
    $ws.Range("B$r").Value = "Vega Modelo de Temuco"
    $ws.Range("C$r").Value = "La Araucanía"
    $ws.Range("D$r").Value2 = $item.Fecha
    $ws.Range("E$r").Value = 9
    $ws.Range("F$r").Value = 100112044
    $ws.Range("G$r").Value = "Perejil"
    $ws.Range("H$r").Value = "Sin especificar"
    $ws.Range("I$r").Value = "Primera"
    $ws.Range("J$r").Value = $item.Volumen
    $ws.Range("K$r").Value = $item.PrecioMin
    $ws.Range("L$r").Value = $item.PrecioMax
    $ws.Range("M$r").Value = $item.PrecioProm
    $ws.Range("N$r").Value = $item.Unidad
    $ws.Range("O$r").Value = $item.Origen
    $ws.Range("P$r").Value = $item.PrecioKg
    $ws.Range("Q$r").Value = 3
    $ws.Range("R$r").Value = "Hortaliza"
}
